$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 2
    3 = 3
    4 = 2
    5 = 2
    6 = 0
    7 = 1
    8 = 5
    9 = 7
    10 = 4
    11 = 4
    12 = 1
    13 = 1
    14 = 2
    15 = 7
    16 = 5
    17 = 4
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
